# Fixed TWOI to SWOI translation: update the LPP t/p/d/BF statistics table
# after re-running the comparisons with the corrected (~8 ms) time window.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Update-Cell($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $ok = $cell.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for row $row col $col ($oldText -> $newText)"
    }
}

Update-Cell 3 2 "3.03" "3.06"   # body 1.t
Update-Cell 3 6 "16.827" "18.027"   # body 1.BF
Update-Cell 4 2 "2.00" "2.01"   # body 2.t
Update-Cell 4 4 ".026" ".025"   # body 2.p
Update-Cell 4 6 " 1.949" " 1.982"   # body 2.BF
Update-Cell 5 2 "-0.74" "-0.75"   # body 3.t
Update-Cell 5 4 ".462" ".458"   # body 3.p
Update-Cell 5 6 " 0.211" " 0.212"   # body 3.BF
Update-Cell 6 2 "1.51" "1.58"   # body 4.t
Update-Cell 6 4 ".073" ".065"   # body 4.p
Update-Cell 6 5 "0.32" "0.34"   # body 4.d
Update-Cell 6 6 " 1.096" " 1.209"   # body 4.BF
Update-Cell 7 2 "2.60" "2.62"   # body 5.t
Update-Cell 7 5 "0.55" "0.56"   # body 5.d
Update-Cell 7 6 " 6.458" " 6.701"   # body 5.BF
Update-Cell 8 2 "0.79" "0.75"   # body 6.t
Update-Cell 8 4 ".439" ".462"   # body 6.p
Update-Cell 8 5 "0.17" "0.16"   # body 6.d
Update-Cell 8 6 " 0.295" " 0.287"   # body 6.BF
Update-Cell 9 2 "2.69" "2.67"   # body 7.t
Update-Cell 9 6 " 7.633" " 7.282"   # body 7.BF
Update-Cell 10 2 "0.72" "0.71"   # body 8.t
Update-Cell 10 4 ".239" ".242"   # body 8.p
Update-Cell 10 6 " 0.424" " 0.420"   # body 8.BF
Update-Cell 11 2 "-2.04" "-2.01"   # body 9.t
Update-Cell 11 4 ".054" ".057"   # body 9.p
Update-Cell 11 6 " 1.256" " 1.213"   # body 9.BF
Update-Cell 12 2 "1.11" "1.03"   # body10.t
Update-Cell 12 4 ".821" ".926"   # body10.p
Update-Cell 12 5 "0.33" "0.31"   # body10.d
Update-Cell 12 6 " 0.488" " 0.456"   # body10.BF
Update-Cell 13 2 "-0.77" "-0.79"   # body11.t
Update-Cell 13 5 "-0.23" "-0.24"   # body11.d
Update-Cell 13 6 " 0.377" " 0.384"   # body11.BF
Update-Cell 14 2 "-1.94" "-1.89"   # body12.t
Update-Cell 14 4 ".176" ".195"   # body12.p
Update-Cell 14 5 "-0.59" "-0.57"   # body12.d
Update-Cell 14 6 " 1.321" " 1.228"   # body12.BF
